# Automatic update of files.
# Swaps the two pairs of duplicate Tretåig hackspett ring-hack observations
# (rows 5/6) and rotates the Garnlav/Järpe/Talltita observations
# (rows 12-15) back into their canonical Id order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 131196446
$ws.Range("Q5").Value = 500268
$ws.Range("R5").Value = 7016174
$ws.Range("AC5").Value = "Ringhack, äldre, på tall."
$ws.Range("AJ5").Value = "tall"
$ws.Range("AK5").Value = "Pinus sylvestris"
$ws.Range("AO5").Value = "Pinus sylvestris"
$ws.Range("A6").Value = 131196443
$ws.Range("Q6").Value = 500146
$ws.Range("R6").Value = 7016293
$ws.Range("AC6").Value = "Ringhack, äldre, på gran."
$ws.Range("AJ6").Value = "gran"
$ws.Range("AK6").Value = "Picea abies"
$ws.Range("AO6").Value = "Picea abies"
$ws.Range("A12").Value = 131196449
$ws.Range("B12").Value = 57064
$ws.Range("E12").Value = 102612
$ws.Range("F12").Value = "Järpe"
$ws.Range("G12").Value = "Tetrastes bonasia"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "2"
$ws.Range("J12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = "födosökande"
$ws.Range("N12").Value = "observerad"
$ws.Range("Q12").Value = 500203
$ws.Range("R12").Value = 7016330
$ws.Range("AC12").Value = "Synobservation av 2 st födosökande järpar."
$ws.Range("AF12").Value = ""
$ws.Range("A13").Value = 131196451
$ws.Range("B13").Value = 79244
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("Q13").Value = 500318
$ws.Range("R13").Value = 7016201
$ws.Range("AC13").Value = ""
$ws.Range("AF13").Value = ""
$ws.Range("A14").Value = 131196447
$ws.Range("B14").Value = 58043
$ws.Range("E14").Value = 103021
$ws.Range("F14").Value = "Talltita"
$ws.Range("G14").Value = "Poecile montanus"
$ws.Range("H14").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "1"
$ws.Range("J14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = "förbiflygande"
$ws.Range("N14").Value = "observerad"
$ws.Range("Q14").Value = 500269
$ws.Range("R14").Value = 7016195
$ws.Range("AC14").Value = "Synobservation av 1 st talltita."
$ws.Range("AF14").Value = ""
$ws.Range("A15").Value = 131196452
$ws.Range("B15").Value = 79244
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("Q15").Value = 500345
$ws.Range("R15").Value = 7016371
$ws.Range("AC15").Value = ""
$ws.Range("AF15").Value = ""
